$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78 (existing rows 78-91 shift down to 79-92)
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with a new price record
$ws.Cells.Item(78, 1).Value = 11
$ws.Cells.Item(78, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(78, 3).Value = "Bíobío"
$ws.Cells.Item(78, 4).Value = 44511
$ws.Cells.Item(78, 5).Value = 8
$ws.Cells.Item(78, 6).Value = "Fruta"
$ws.Cells.Item(78, 7).Value = 100102
$ws.Cells.Item(78, 8).Value = "Cítricos"
$ws.Cells.Item(78, 9).Value = 100102004
$ws.Cells.Item(78, 10).Value = "Mandarina"
$ws.Cells.Item(78, 11).Value = "Murcott"
$ws.Cells.Item(78, 12).Value = "Primera"
$ws.Cells.Item(78, 13).Value = 350
$ws.Cells.Item(78, 14).Value = 7000
$ws.Cells.Item(78, 15).Value = 7500
$ws.Cells.Item(78, 16).Value = 7286
$ws.Cells.Item(78, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(78, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(78, 19).Value = 405
$ws.Cells.Item(78, 20).Value = 18
